# "started chapter 25 (i guess.. or 24..)"
# Add a new glossary entry below the existing rows (row 29 -> new row 30):
#   Nr. = 26, Wort = "Kernel/Filter", Anmerkung = "Siehe concolutional Layer"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$newRow = 30

# Write the "Anmerkung" (column D) first and the glossary word (column B) second so
# that the new shared-string entries are appended in the same order as the source
# workbook ("Siehe concolutional Layer" before "Kernel/Filter").
$ws.Cells.Item($newRow, 4).Value = "Siehe concolutional Layer"
$ws.Cells.Item($newRow, 1).Value = 26
$ws.Cells.Item($newRow, 2).Value = "Kernel/Filter"

# Nudge the column widths for A/B to the values Excel settled on after re-opening
# the workbook (closest values reachable through this engine's width quantisation).
$ws.Columns.Item(1).ColumnWidth = 19.8
$ws.Columns.Item(2).ColumnWidth = 23.5

$ws.Range("B30").Select()
